$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records need to be inserted right before the current row 298,
# pushing the existing rows 298-310 down to become rows 300-312.
$ws.Rows.Item(298).Insert()
$ws.Rows.Item(298).Insert()

# --- New row 298 ---
$ws.Cells.Item(298,1).Value = 7
$ws.Cells.Item(298,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(298,3).Value = "Ñuble"
$ws.Cells.Item(298,4).Value = 44568
$ws.Cells.Item(298,4).NumberFormat = $ws.Cells.Item(297,4).NumberFormat
$ws.Cells.Item(298,5).Value = 16
$ws.Cells.Item(298,6).Value = 100114014
$ws.Cells.Item(298,7).Value = "Betarraga"
$ws.Cells.Item(298,8).Value = "Sin especificar"
$ws.Cells.Item(298,9).Value = "Primera"
$ws.Cells.Item(298,10).Value = 600
$ws.Cells.Item(298,11).Value = 500
$ws.Cells.Item(298,12).Value = 600
$ws.Cells.Item(298,13).Value = 550
$ws.Cells.Item(298,14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(298,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(298,16).Value = 110
$ws.Cells.Item(298,17).Value = 5
$ws.Cells.Item(298,18).Value = "Hortaliza"

# --- New row 299 ---
$ws.Cells.Item(299,1).Value = 7
$ws.Cells.Item(299,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(299,3).Value = "Ñuble"
$ws.Cells.Item(299,4).Value = 44568
$ws.Cells.Item(299,4).NumberFormat = $ws.Cells.Item(297,4).NumberFormat
$ws.Cells.Item(299,5).Value = 16
$ws.Cells.Item(299,6).Value = 100114014
$ws.Cells.Item(299,7).Value = "Betarraga"
$ws.Cells.Item(299,8).Value = "Sin especificar"
$ws.Cells.Item(299,9).Value = "Primera"
$ws.Cells.Item(299,10).Value = 600
$ws.Cells.Item(299,11).Value = 500
$ws.Cells.Item(299,12).Value = 600
$ws.Cells.Item(299,13).Value = 550
$ws.Cells.Item(299,14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(299,15).Value = "Región del Maule"
$ws.Cells.Item(299,16).Value = 110
$ws.Cells.Item(299,17).Value = 5
$ws.Cells.Item(299,18).Value = "Hortaliza"
